# Apply permutation of rows 2-23 (columns D, L, M, N, O, P, R, S) per commit
# "Fruta / hortaliza, semanal" — the weekly refresh reshuffled which market
# observation lands on which spreadsheet row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that move, keyed by row,
# so every write below reads from this fixed snapshot instead of any
# already-updated cell (rows are permuted, so some values must survive
# being overwritten elsewhere on the sheet).
$D = @{}
$L = @{}
$M = @{}
$N = @{}
$O = @{}
$P = @{}
$R = @{}
$S = @{}
for ($r = 2; $r -le 23; $r++) {
    $D[$r] = $ws.Cells.Item($r, 4).Value()
    $L[$r] = $ws.Cells.Item($r, 12).Value()
    $M[$r] = $ws.Cells.Item($r, 13).Value()
    $N[$r] = $ws.Cells.Item($r, 14).Value()
    $O[$r] = $ws.Cells.Item($r, 15).Value()
    $P[$r] = $ws.Cells.Item($r, 16).Value()
    $R[$r] = $ws.Cells.Item($r, 18).Value()
    $S[$r] = $ws.Cells.Item($r, 19).Value()
}

# destination row -> source row (the row whose pre-edit data now belongs here)
$map = @{
    2 = 16
    3 = 17
    4 = 8
    5 = 11
    6 = 12
    7 = 10
    8 = 23
    9 = 5
    10 = 6
    11 = 4
    12 = 18
    13 = 19
    14 = 9
    15 = 20
    16 = 2
    17 = 3
    18 = 13
    19 = 14
    20 = 15
    21 = 7
    22 = 21
    23 = 22
}

foreach ($destRow in ($map.Keys | Sort-Object)) {
    $srcRow = $map[$destRow]
    $ws.Cells.Item($destRow, 4).Value = $D[$srcRow]
    $ws.Cells.Item($destRow, 12).Value = $L[$srcRow]
    $ws.Cells.Item($destRow, 13).Value = $M[$srcRow]
    $ws.Cells.Item($destRow, 14).Value = $N[$srcRow]
    $ws.Cells.Item($destRow, 15).Value = $O[$srcRow]
    $ws.Cells.Item($destRow, 16).Value = $P[$srcRow]
    $ws.Cells.Item($destRow, 18).Value = $R[$srcRow]
    $ws.Cells.Item($destRow, 19).Value = $S[$srcRow]
}
